# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-06-05 Thursday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-06-06 Friday", 2)

# Update the multiplication problems in the table. Cells are addressed
# positionally (row/column) rather than by old text, because one of the
# new values ("775×2=") collides with an old value that appears later in
# the table and must itself be replaced with a different new value.
$t = $d.Tables.Item(1)

$rows = @(1, 5, 10, 15, 20)
$newValues = @(
    @("981×7=", "926×7=", "472×4=", "875×3=", "476×7="),
    @("693×8=", "338×3=", "986×7=", "702×2=", "448×4="),
    @("775×2=", "177×2=", "285×9=", "684×8=", "802×5="),
    @("206×8=", "761×2=", "498×8=", "882×5=", "427×7="),
    @("145×8=", "679×6=", "121×7=", "829×8=", "540×7=")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$i][$c - 1]
    }
}
